$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

$ws.Range("B2").Value = 'Component: Multi-Functional Tool Application'
$ws.Range("E3").Value = 'MFP: Any'
$ws.Range("E6").Value = '1. Copy MultiFunctionalTool_For_Desktop.zip from specified tec-share location<br>2. Extract contents to preferred location<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range("F6").Value = 'Application launches successfully with all features accessible'
$ws.Range("E7").Value = '1. Copy MultiFunctionalTool_For_Laptop.zip from specified tec-share location<br>2. Extract contents to preferred location<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range("F7").Value = 'Application launches successfully with all features accessible'
$ws.Range("C8").Value = 'Application is installed and running'
$ws.Range("D8").Value = 'Verify Network Packet Capture start functionality'
$ws.Range("E8").Value = '1. Navigate to Network Packet Capture section<br>2. Click Start button'
$ws.Range("F8").Value = 'Packet capture begins successfully'
$ws.Range("C9").Value = 'Network Packet Capture is running'
$ws.Range("D9").Value = 'Verify Network Packet Capture stop functionality'
$ws.Range("E9").Value = '1. With packet capture running, click Stop button'
$ws.Range("F9").Value = 'Packet capture stops and generates .pcap file'
$ws.Range("C10").Value = 'Network Packet Capture has been completed'
$ws.Range("D10").Value = 'Verify .pcap file storage'
$ws.Range("E10").Value = '1. Complete a packet capture session<br>2. Check MFP''s Shared Folder'
$ws.Range("F10").Value = '.pcap file is automatically copied to MFP''s Shared Folder and folder opens automatically'
$ws.Range("C11").Value = 'Application is installed and running'
$ws.Range("E11").Value = '1. Navigate to Memory Leak Check section<br>2. Select a protocol<br>3. Initiate memory leak check'
$ws.Range("F11").Value = 'Memory Leak Comparison Table is displayed with accurate information'
$ws.Range("C12").Value = 'Application is installed and running'
$ws.Range("E12").Value = '1. Navigate to Debug Log Collection section<br>2. Click Run button'
$ws.Range("F12").Value = 'Script executes and begins collecting logs'
$ws.Range("C13").Value = 'Debug Log Collection has been completed'
$ws.Range("D13").Value = 'Verify debug log storage'
$ws.Range("E13").Value = '1. Complete a debug log collection<br>2. Check MFP''s Shared Folder'
$ws.Range("F13").Value = 'Debug logs are copied to MFP''s Shared Folder and folder opens automatically'
$ws.Range("C14").Value = 'Debug Log Collection folder is empty on first attempt'
$ws.Range("D14").Value = 'Verify debug log collection retry functionality'
$ws.Range("E14").Value = '1. Navigate to Debug Log Collection section<br>2. Click Run button<br>3. If folder is empty, run operation again'
$ws.Range("F14").Value = 'Debug logs are successfully collected on second attempt'
$ws.Range("C15").Value = 'Application is installed and running'
$ws.Range("D15").Value = 'Verify Diagnostic Code Details - ECC'
$ws.Range("E15").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select ECC option'
$ws.Range("F15").Value = 'Relevant ECC diagnostic details are displayed'
$ws.Range("C16").Value = 'Application is installed and running'
$ws.Range("D16").Value = 'Verify Diagnostic Code Details - Network Protocols'
$ws.Range("E16").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select Network Protocols option'
$ws.Range("F16").Value = 'Relevant Network Protocols diagnostic details are displayed'
$ws.Range("C17").Value = 'Application is installed and running'
$ws.Range("D17").Value = 'Verify Diagnostic Code Details - High Security Mode'
$ws.Range("E17").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select High Security Mode option'
$ws.Range("F17").Value = 'Relevant High Security Mode diagnostic details are displayed'
$ws.Range("C18").Value = 'Application is installed and running'
$ws.Range("D18").Value = 'Verify Diagnostic Code Details - Common codes'
$ws.Range("E18").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select a commonly used diagnostic code'
$ws.Range("F18").Value = 'Relevant job-specific details for the selected code are displayed'
$ws.Range("C19").Value = 'Application is installed and running'
$ws.Range("D19").Value = 'Verify 08 Diagnostic Code Value - Get functionality'
$ws.Range("E19").Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Click Get button'
$ws.Range("F19").Value = 'Current value of the selected diagnostic code is displayed'
$ws.Range("C20").Value = 'Application is installed and running'
$ws.Range("D20").Value = 'Verify 08 Diagnostic Code Value - Set functionality'
$ws.Range("E20").Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Enter a new value<br>4. Click Set button'
$ws.Range("F20").Value = 'Diagnostic code value is updated successfully'
$ws.Range("C21").Value = 'Application is installed and running'
$ws.Range("D21").Value = 'Verify Protocol Configuration - Get functionality'
$ws.Range("E21").Value = '1. Navigate to Protocol Configuration section<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Click Get button'
$ws.Range("F21").Value = 'Current value of the selected protocol is displayed'
$ws.Range("C22").Value = 'Application is installed and running'
$ws.Range("D22").Value = 'Verify Protocol Configuration - Set functionality'
$ws.Range("E22").Value = '1. Navigate to Protocol Configuration section<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Enter a new value<br>5. Click Set button'
$ws.Range("F22").Value = 'Message indicating "Set protocol values operation still needs to be implemented" is displayed'
$ws.Range("H22").Value = 'Feature not yet implemented'
$ws.Range("C23").Value = 'Application is running with multiple features'
$ws.Range("D23").Value = 'Verify GUI responsiveness'
$ws.Range("E23").Value = '1. Launch application<br>2. Navigate between different features rapidly<br>3. Perform operations in quick succession'
$ws.Range("F23").Value = 'GUI remains responsive with no freezing or significant delays'
$ws.Range("H23").Value = ""
$ws.Range("B24").Value = 'TC019'
$ws.Range("C24").Value = 'Application is running on a system with limited resources'
$ws.Range("D24").Value = 'Verify performance under resource constraints'
$ws.Range("E24").Value = '1. Launch application on a system with minimal RAM/CPU<br>2. Perform all main functions sequentially'
$ws.Range("F24").Value = 'Application performs all functions without crashing or excessive resource usage'
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = ""
$ws.Range("B25").Value = 'TC020'
$ws.Range("C25").Value = 'Application is installed and running'
$ws.Range("D25").Value = 'Verify time efficiency for packet capture'
$ws.Range("E25").Value = '1. Measure time to complete packet capture manually<br>2. Measure time to complete same packet capture using the application'
$ws.Range("F25").Value = 'Application reduces testing time by approximately 80% compared to manual method'
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = ""
$ws.Range("B26").Value = 'TC021'
$ws.Range("C26").Value = 'Application is installed and running'
$ws.Range("D26").Value = 'Verify usability for new users'
$ws.Range("E26").Value = '1. Provide application to a user unfamiliar with it<br>2. Ask them to perform basic operations without instructions<br>3. Observe and record their experience'
$ws.Range("F26").Value = 'New user can navigate and use basic features without significant confusion'
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = ""
$ws.Range("B27").Value = 'TC022'
$ws.Range("C27").Value = 'Application is installed on different Windows OS versions'
$ws.Range("D27").Value = 'Verify OS compatibility'
$ws.Range("E27").Value = '1. Install and run application on Windows 7, 8, 10, and 11<br>2. Test core functionality on each OS'
$ws.Range("F27").Value = 'Application functions correctly on all tested Windows versions'
$ws.Range("G27").Value = ""
$ws.Range("H27").Value = ""
$ws.Range("B28").Value = 'TC023'
$ws.Range("C28").Value = 'Application is running'
$ws.Range("D28").Value = 'Verify error handling for invalid inputs'
$ws.Range("E28").Value = '1. Enter invalid values in input fields<br>2. Submit the invalid data'
$ws.Range("F28").Value = 'Application displays appropriate error messages without crashing'
$ws.Range("G28").Value = ""
$ws.Range("H28").Value = ""
$ws.Range("B29").Value = 'TC024'
$ws.Range("C29").Value = 'Application is running'
$ws.Range("D29").Value = 'Verify simultaneous operations'
$ws.Range("E29").Value = '1. Attempt to run multiple operations simultaneously<br>2. Observe application behavior'
$ws.Range("F29").Value = 'Application either handles concurrent operations correctly or provides clear feedback about limitations'
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = ""
$ws.Range("B30").Value = 'TC025'
$ws.Range("C30").Value = 'Application is running'
$ws.Range("D30").Value = 'Verify data persistence'
$ws.Range("E30").Value = '1. Configure settings in the application<br>2. Close and reopen the application'
$ws.Range("F30").Value = 'Previously configured settings are retained'
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = ""

# Extend used range to column L to match target dimension (A1:L111),
# reusing the existing hidden-column style (46) so no new style is introduced.
$ws.Range("L111").Value = "x"
$ws.Range("L111").Value = ""

